$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 111, pushing the existing rows 111..143 down to 112..144
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with the new data point
$ws.Cells.Item(111,1).Value  = 7
$ws.Cells.Item(111,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(111,3).Value  = "Ñuble"
$ws.Cells.Item(111,4).Value  = 44508
$ws.Cells.Item(111,5).Value  = 16
$ws.Cells.Item(111,6).Value  = 100112006
$ws.Cells.Item(111,7).Value  = "Repollo"
$ws.Cells.Item(111,8).Value  = "Crespo record"
$ws.Cells.Item(111,9).Value  = "Primera"
$ws.Cells.Item(111,10).Value = 400
$ws.Cells.Item(111,11).Value = 600
$ws.Cells.Item(111,12).Value = 700
$ws.Cells.Item(111,13).Value = 650
$ws.Cells.Item(111,14).Value = "$/unidad"
$ws.Cells.Item(111,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(111,16).Value = 650
$ws.Cells.Item(111,17).Value = 1
$ws.Cells.Item(111,18).Value = "Hortaliza"
